$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level value corrections in rows 1-25 (no row shifting involved) ---
$ws.Range("D5").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("D11").Value = -15.5
$ws.Range("C19").Value = 13.2
$ws.Range("D19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("D25").Value = -15.5

# --- Remove the rows for "RM 232" and "SC 92" entirely (data cleanup) ---
# Delete higher row index first so the lower one's index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Cell-level value corrections in the rows following the deletion (now renumbered 26-33) ---
$ws.Range("B26").ClearContents()
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E28").Value = -5.9
$ws.Range("B29").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
